$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cryptocurrency price/volume data refresh (GitHub Actions scheduled update).
# Column D ("Price") values are kept as literal text, matching the source data feed,
# so we force the cell number format to Text before assigning them - this prevents
# Excel from re-interpreting strings such as "1.001" or "0.08480" as numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '23.918.48'
$ws.Range('E2').Value = '  -2.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.653.68'
$ws.Range('E3').Value = '  -1.02%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.44'
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3891'
$ws.Range('E7').Value = '  -2.00%  '
$ws.Range('E8').Value = '  -2.88%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '51.45'
$ws.Range('E9').Value = '  -1.71%  '
$ws.Range('E10').Value = '  -3.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9997'
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08480'
$ws.Range('E13').Value = '  -2.29%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.048'
$ws.Range('E14').Value = '  -3.60%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.064'
$ws.Range('E15').Value = '  +1.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001316'
$ws.Range('E16').Value = '  -1.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.654.62'
$ws.Range('E17').Value = '  -0.63%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '94.17'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06997'
$ws.Range('E19').Value = '  -0.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.61'
$ws.Range('E20').Value = '  -5.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.985'
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.71'
$ws.Range('E23').Value = '  -0.52%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '23.904.68'
$ws.Range('E24').Value = '  -2.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.433'
$ws.Range('E25').Value = '  -1.84%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.961'
$ws.Range('E26').Value = '  -3.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.04'
$ws.Range('E27').Value = '  -2.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '153.89'
$ws.Range('E28').Value = '  -2.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.427'
$ws.Range('E29').Value = '  -0.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '137.63'
$ws.Range('E30').Value = '  -3.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.802'
$ws.Range('E31').Value = '  -2.77%  '
$ws.Range('E32').Value = '  -2.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.844.00'
$ws.Range('E33').Value = '  -0.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08160'
$ws.Range('E34').Value = '  -1.65%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.007'
$ws.Range('E35').Value = '  -5.59%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02915'
$ws.Range('E36').Value = '  -6.24%  '
$ws.Range('E37').Value = '  -3.78%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '10.79'
$ws.Range('E38').Value = '  -3.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2675'
$ws.Range('E39').Value = '  -3.43%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.09131'
$ws.Range('E40').Value = '  -1.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '13.57'
$ws.Range('E41').Value = '  -1.48%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7559'
$ws.Range('E42').Value = '  -2.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.424'
$ws.Range('E43').Value = '  -1.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.51'
$ws.Range('E44').Value = '  -0.64%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6932'
$ws.Range('E45').Value = '  -2.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.453'
$ws.Range('E46').Value = '  -3.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.106'
$ws.Range('E47').Value = '  -0.61%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9992'
$ws.Range('E48').Value = '  -0.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08281'
$ws.Range('E49').Value = '  -1.98%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '133.24'
$ws.Range('E50').Value = '  -2.86%  '
$ws.Range('E51').Value = '  -3.38%  '
